$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BenchmarkResults")
Write-Host $ws.Name
